$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 80780
$ws.Range("J136").Value = 80780
$ws.Range("L136").Value = 80780
$ws.Range("N136").Value = -90980
$ws.Range("H137").Value = 1369.2593
$ws.Range("I137").Value = 1198.3334
$ws.Range("J137").Value = 1711.1111
$ws.Range("K137").Value = 3595.0002
$ws.Range("L137").Value = 5133.3333
$ws.Range("M137").Value = -1045.0002
$ws.Range("N137").Value = -10233.3333
$ws.Range("H138").Value = 2720
$ws.Range("I138").Value = 2602.7334
$ws.Range("J138").Value = 2782.8215
$ws.Range("K138").Value = 7808.2002
$ws.Range("L138").Value = 8348.4645
$ws.Range("M138").Value = -2668.2002
$ws.Range("N138").Value = -18628.4645

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6056.5454
$ws.Range("I28").Value = 3862.2
$ws.Range("K28").Value = 3862.2
$ws.Range("M28").Value = -3670.2
$ws.Range("H31").Value = 4533.3335
$ws.Range("I31").Value = 4533.3335
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 4533.3335
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -4239.3335
$ws.Range("N31").ClearContents()
$ws.Range("H32").Value = 24497.637
$ws.Range("I32").Value = 4829.8115
$ws.Range("J32").Value = 104681.84
$ws.Range("K32").Value = 4829.8115
$ws.Range("L32").Value = 104681.84
$ws.Range("M32").Value = -4542.8115
$ws.Range("N32").Value = -105255.84
$ws.Range("H61").Value = 4715.091
$ws.Range("I61").Value = 2838.2856
$ws.Range("J61").Value = 7999.5
$ws.Range("K61").Value = 2838.2856
$ws.Range("L61").Value = 7999.5
$ws.Range("M61").Value = -2626.2856
$ws.Range("N61").Value = -8423.5
$ws.Range("H74").Value = 1302.9286
$ws.Range("I74").Value = 865.0645
$ws.Range("K74").Value = 865.0645
$ws.Range("M74").Value = 8.935500000000047
$ws.Range("H77").Value = 1302.9286
$ws.Range("I77").Value = 865.0645
$ws.Range("K77").Value = 4325.3225
$ws.Range("M77").Value = 42.67749999999978
$ws.Range("H99").Value = 6056.5454
$ws.Range("I99").Value = 3862.2
$ws.Range("K99").Value = 3862.2
$ws.Range("M99").Value = -867.1999999999998
$ws.Range("H118").Value = 32222.285
$ws.Range("J118").Value = 32222.285
$ws.Range("L118").Value = 32222.285
$ws.Range("N118").Value = -35536.285
$ws.Range("H122").Value = 2843.2222
$ws.Range("I122").Value = 2198.3333
$ws.Range("J122").Value = 3165.6667
$ws.Range("K122").Value = 6594.999899999999
$ws.Range("L122").Value = 9497.000100000001
$ws.Range("M122").Value = -4144.999899999999
$ws.Range("N122").Value = -14397.0001
$ws.Range("H136").Value = 4715.091
$ws.Range("I136").Value = 2838.2856
$ws.Range("J136").Value = 7999.5
$ws.Range("K136").Value = 8514.856800000001
$ws.Range("L136").Value = 23998.5
$ws.Range("M136").Value = -5964.856800000001
$ws.Range("N136").Value = -29098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 14180.667
$ws.Range("J21").Value = 14180.667
$ws.Range("L21").Value = 14180.667
$ws.Range("N21").Value = -14652.667
$ws.Range("H27").Value = 22116.666
$ws.Range("J27").Value = 22116.666
$ws.Range("L27").Value = 22116.666
$ws.Range("N27").Value = -22500.666
$ws.Range("H63").Value = 35757
$ws.Range("J63").Value = 35757
$ws.Range("L63").Value = 35757
$ws.Range("N63").Value = -37129
$ws.Range("H66").Value = 35757
$ws.Range("J66").Value = 35757
$ws.Range("L66").Value = 107271
$ws.Range("N66").Value = -114135
$ws.Range("H97").Value = 61326.43
$ws.Range("I97").Value = 70547.5
$ws.Range("J97").Value = 6000
$ws.Range("K97").Value = 70547.5
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -69556.5
$ws.Range("N97").Value = -7982
$ws.Range("H98").Value = 25000
$ws.Range("J98").Value = 25000
$ws.Range("L98").Value = 25000
$ws.Range("N98").Value = -30990
$ws.Range("H105").Value = 3233.3333
$ws.Range("I105").Value = 2388.889
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 2388.889
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -641.8890000000001
$ws.Range("N105").Value = -7994
$ws.Range("H111").Value = 34998
$ws.Range("J111").Value = 34998
$ws.Range("L111").Value = 34998
$ws.Range("N111").Value = -43178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2453.2856
$ws.Range("I31").Value = 2050.4583
$ws.Range("J31").Value = 2701.1794
$ws.Range("K31").Value = 2050.4583
$ws.Range("L31").Value = 2701.1794
$ws.Range("M31").Value = -1755.4583
$ws.Range("N31").Value = -3291.1794
$ws.Range("H34").Value = 2453.2856
$ws.Range("I34").Value = 2050.4583
$ws.Range("J34").Value = 2701.1794
$ws.Range("K34").Value = 2050.4583
$ws.Range("L34").Value = 2701.1794
$ws.Range("M34").Value = -1848.4583
$ws.Range("N34").Value = -3105.1794
$ws.Range("H50").Value = 13194
$ws.Range("J50").Value = 13194
$ws.Range("L50").Value = 13194
$ws.Range("N50").Value = -14444
$ws.Range("H51").Value = 18451.6
$ws.Range("J51").Value = 18451.6
$ws.Range("L51").Value = 18451.6
$ws.Range("N51").Value = -19923.6
$ws.Range("H59").Value = 250021680
$ws.Range("J59").Value = 250021680
$ws.Range("L59").Value = 250021680
$ws.Range("N59").Value = -250023970
$ws.Range("H60").Value = 22990
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 22990
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 22990
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -24012
$ws.Range("H61").Value = 18451.6
$ws.Range("J61").Value = 18451.6
$ws.Range("L61").Value = 18451.6
$ws.Range("N61").Value = -19147.6
$ws.Range("H87").Value = 26920
$ws.Range("J87").Value = 26920
$ws.Range("L87").Value = 26920
$ws.Range("N87").Value = -29292
$ws.Range("H90").Value = 26920
$ws.Range("J90").Value = 26920
$ws.Range("L90").Value = 80760
$ws.Range("N90").Value = -92616
$ws.Range("H97").Value = 9780
$ws.Range("J97").Value = 9780
$ws.Range("L97").Value = 9780
$ws.Range("N97").Value = -11762
$ws.Range("H105").Value = 4358.727
$ws.Range("I105").Value = 3798.5386
$ws.Range("K105").Value = 3798.5386
$ws.Range("M105").Value = -2051.5386
$ws.Range("H107").Value = 1454.9474
$ws.Range("I107").Value = 1320.7858
$ws.Range("J107").Value = 1830.6
$ws.Range("K107").Value = 1320.7858
$ws.Range("L107").Value = 1830.6
$ws.Range("M107").Value = 599.2141999999999
$ws.Range("N107").Value = -5670.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 530.7083
$ws.Range("I5").Value = 418.44446
$ws.Range("J5").Value = 867.5
$ws.Range("K5").Value = 1255.33338
$ws.Range("L5").Value = 2602.5
$ws.Range("M5").Value = -1143.33338
$ws.Range("N5").Value = -2826.5
$ws.Range("H68").Value = 2234.907
$ws.Range("I68").Value = 1940.6364
$ws.Range("J68").Value = 2543.1904
$ws.Range("K68").Value = 5821.9092
$ws.Range("L68").Value = 7629.5712
$ws.Range("M68").Value = -5010.9092
$ws.Range("N68").Value = -9251.5712
$ws.Range("H71").Value = 2234.907
$ws.Range("I71").Value = 1940.6364
$ws.Range("J71").Value = 2543.1904
$ws.Range("K71").Value = 17465.7276
$ws.Range("L71").Value = 22888.7136
$ws.Range("M71").Value = -13409.7276
$ws.Range("N71").Value = -31000.7136
$ws.Range("H107").Value = 876.9298
$ws.Range("I107").Value = 520.53845
$ws.Range("J107").Value = 1649.1111
$ws.Range("K107").Value = 1561.61535
$ws.Range("L107").Value = 4947.3333
$ws.Range("M107").Value = 358.38465
$ws.Range("N107").Value = -8787.3333
$ws.Range("H131").Value = 930.3333
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 996.53845
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2989.61535
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -13069.61535
$ws.Range("H135").Value = 530.7083
$ws.Range("I135").Value = 418.44446
$ws.Range("J135").Value = 867.5
$ws.Range("K135").Value = 3766.00014
$ws.Range("L135").Value = 7807.5
$ws.Range("M135").Value = -1231.00014
$ws.Range("N135").Value = -12877.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27514.596
$ws.Range("I70").Value = 47407.08
$ws.Range("K70").Value = 47407.08
$ws.Range("M70").Value = -47137.08
$ws.Range("H73").Value = 27514.596
$ws.Range("I73").Value = 47407.08
$ws.Range("K73").Value = 47407.08
$ws.Range("M73").Value = -46471.08
$ws.Range("H99").Value = 5248.4
$ws.Range("I99").Value = 2810.5
$ws.Range("J99").Value = 15000
$ws.Range("K99").Value = 2810.5
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -564.5
$ws.Range("N99").Value = -19492
$ws.Range("H120").Value = 36444
$ws.Range("J120").Value = 36444
$ws.Range("L120").Value = 36444
$ws.Range("N120").Value = -46120
$ws.Range("H132").Value = 2167.0256
$ws.Range("I132").Value = 1633.6
$ws.Range("J132").Value = 3945.111
$ws.Range("K132").Value = 4900.799999999999
$ws.Range("L132").Value = 11835.333
$ws.Range("M132").Value = -2370.799999999999
$ws.Range("N132").Value = -16895.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 142860130
$ws.Range("I16").Value = 142860130
$ws.Range("K16").Value = 142860130
$ws.Range("M16").Value = -142859960
$ws.Range("H100").Value = 2644.7778
$ws.Range("I100").Value = 2257.5715
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 2257.5715
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -1716.5715
$ws.Range("N100").Value = -5082
$ws.Range("H124").Value = 33000
$ws.Range("J124").Value = 33000
$ws.Range("L124").Value = 33000
$ws.Range("N124").Value = -42820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 100000000
$ws.Range("J15").Value = 100000000
$ws.Range("L15").Value = 100000000
$ws.Range("N15").Value = -100000576

Write-Host "Applied all changes"